# "Maly upravy v ramci citelnosti" (small readability edits): reorder the rows of the
# Frame0 table and insert two new course rows (seminarici_bez_seminare.xlsx).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Stage the original nazev/zkratka columns (A2:B32) off to one side so that cell
#    TYPES (text vs number -- e.g. the zkratka "0153" must stay text) survive the
#    reshuffle: Copy/PasteSpecial preserves the original string typing, whereas just
#    assigning Value2 on a numeric-looking string gets silently coerced into a number.
$ws.Range("A2:B32").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = $false

# 2) Rebuild rows 2..34 in their final order, pulling nazev/zkratka back in from the
#    staged copy, and writing the seminariciUcitIdno number directly.
$ws.Range("E2:F2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(2, 3).Value2 = 14

$ws.Range("E3:F3").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(3, 3).Value2 = 14

$ws.Range("E4:F4").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(4, 3).Value2 = 302

$ws.Range("E5:F5").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(5, 3).Value2 = 306

$ws.Range("E7:F7").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(6, 3).Value2 = 313

$ws.Range("E6:F6").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(7, 3).Value2 = 313

$ws.Range("E8:F8").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(8, 3).Value2 = 612

$ws.Range("E9:F9").Copy() | Out-Null
$ws.Range("A9").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(9, 3).Value2 = 612

$ws.Range("E10:F10").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(10, 3).Value2 = 612

$ws.Range("E11:F11").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(11, 3).Value2 = 1609

$ws.Range("E12:F12").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(12, 3).Value2 = 1609

$ws.Range("E13:F13").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(13, 3).Value2 = 2527

$ws.Range("E14:F14").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(14, 3).Value2 = 3457

$ws.Range("E15:F15").Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(15, 3).Value2 = 3457

$ws.Range("E16:F16").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(16, 3).Value2 = 3606

$ws.Range("E17:F17").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(17, 3).Value2 = 3606

$ws.Range("E18:F18").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(18, 3).Value2 = 4190

$ws.Range("E19:F19").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(19, 3).Value2 = 4221

$ws.Range("E20:F20").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(20, 3).Value2 = 4625

$ws.Range("E22:F22").Copy() | Out-Null
$ws.Range("A21").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(21, 3).Value2 = 4746

$ws.Range("E21:F21").Copy() | Out-Null
$ws.Range("A22").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(22, 3).Value2 = 4746

$ws.Cells.Item(23, 1).Value2 = "Scientific inquiry and reasoning"
$ws.Cells.Item(23, 2).Value2 = "PD101"
$ws.Cells.Item(23, 3).Value2 = 4747

$ws.Range("E23:F23").Copy() | Out-Null
$ws.Range("A24").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(24, 3).Value2 = 4991

$ws.Range("E24:F24").Copy() | Out-Null
$ws.Range("A25").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(25, 3).Value2 = 4991

$ws.Cells.Item(26, 1).Value2 = "Letní geografická škola"
$ws.Range("Z1").Formula = '="0158"'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("B26").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").Clear() | Out-Null
$ws.Cells.Item(26, 3).Value2 = 5886

$ws.Range("E26:F26").Copy() | Out-Null
$ws.Range("A27").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(27, 3).Value2 = 8021

$ws.Range("E27:F27").Copy() | Out-Null
$ws.Range("A28").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(28, 3).Value2 = 8021

$ws.Range("E25:F25").Copy() | Out-Null
$ws.Range("A29").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(29, 3).Value2 = 8021

$ws.Range("E29:F29").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(30, 3).Value2 = 8093

$ws.Range("E28:F28").Copy() | Out-Null
$ws.Range("A31").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(31, 3).Value2 = 8093

$ws.Range("E31:F31").Copy() | Out-Null
$ws.Range("A32").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(32, 3).Value2 = 8514

$ws.Range("E32:F32").Copy() | Out-Null
$ws.Range("A33").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(33, 3).Value2 = 8514

$ws.Range("E30:F30").Copy() | Out-Null
$ws.Range("A34").PasteSpecial(-4104) | Out-Null
$ws.Cells.Item(34, 3).Value2 = 8514

$excel.CutCopyMode = $false

# 3) New rows (33/34) need the same body formatting (vertical-centered text, and the
#    custom thousands number format on column C) as the rest of the table; copy that
#    formatting down from the last original body row.
$ws.Range("A32:C32").Copy() | Out-Null
$ws.Range("A33:C34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 4) Clear the staging columns again.
$ws.Range("E2:F32").Clear() | Out-Null

# 5) Grow the table (ListObject) + autofilter to the two extra rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C34"))

